$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous (small) table before laying out the new, larger one.
$ws.Cells.Clear()

# Header row: company names
$ws.Range("C2").Value = "捷揚"
$ws.Range("D2").Value = "新光"
$ws.Range("E2").Value = "中興"

# Row 3: 月費 (monthly fee)
$ws.Range("B3").Value = "月費"
$ws.Range("C3").Value = 2000
$ws.Range("D3").Value = 1600

# Row 4: 安裝費 (installation fee)
$ws.Range("B4").Value = "安裝費"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1000

# Row 5: 最少簽約期限 (minimum contract period)
$ws.Range("B5").Value = "最少簽約期限"
$ws.Range("C5").Value = "一年"
$ws.Range("D5").Value = "繳一年算13個月"

# Row 6: 保全範圍 (coverage scope)
$ws.Range("B6").Value = "保全範圍"
$ws.Range("C6").Value = "鐵捲門 一般門 窗戶"
$ws.Range("D6").Value = "鐵捲門 一般門 窗戶"

# Row 7: 有無監視器 (camera available?)
$ws.Range("B7").Value = "有無監視器"
$ws.Range("C7").Value = "無"
$ws.Range("D7").Value = "無"

# Row 8: 保險內容 (insurance content)
$ws.Range("B8").Value = "保險內容"
$ws.Range("C8").Value = "竊盜險40萬"
$ws.Range("D8").Value = "最高理賠200倍"

# Row 9: 使用網路 (network used)
$ws.Range("B9").Value = "使用網路"
$ws.Range("C9").Value = "固定基地台，可支援sim卡"
$ws.Range("D9").Value = "固定基地台"

# Row 10: 配合鎖匠 (locksmith cooperation)
$ws.Range("B10").Value = "配合鎖匠"
$ws.Range("C10").Value = "無"
$ws.Range("D10").Value = "無"

# Row 11: 可否假日安裝 (installable on holidays?)
$ws.Range("B11").Value = "可否假日安裝"
$ws.Range("C11").Value = "可以星期日(要安排"
$ws.Range("D11").Value = "可以星期日"

# Row 12: 安裝前幾天聯絡 (days before installation contact)
$ws.Range("B12").Value = "安裝前幾天聯絡"
$ws.Range("C12").Value = "三天之前(要發包"
$ws.Range("D12").Value = "當天聯絡就可以"

# Row 13: 聯絡電話 (contact numbers)
$ws.Range("B13").Value = "聯絡電話"
$ws.Range("C13").Value = "'0905 300 313許先生"
$ws.Range("D13").Value = "0982 819 655林先生"

# Row 14: extra phone number for 中興 in column E
$ws.Range("E14").Value = "835 7811"

# Resize the columns to fit the new (wider) content, matching the
# widths Excel's own AutoFit produced for this data.
$ws.Columns.Item(1).ColumnWidth = 8.285714285714286
$ws.Columns.Item(2).ColumnWidth = 26.428571428571427
$ws.Columns.Item(3).ColumnWidth = 25.142857142857142
$ws.Columns.Item(4).ColumnWidth = 25.142857142857142
$ws.Columns.Item(5).ColumnWidth = 28.0

# Restore the selection to match the saved view state.
$ws.Range("E18").Select()
